$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "r775"
$ws.Range("B15").Value = "cameron"
$ws.Range("C15").Value = "is this the way?"
$ws.Range("D15").Value = "2025-10-01 14:54:24"
